function Set-RowValues($ws, $row, $a, $b, $c, $d, $e) {
    $ws.Cells.Item($row, 1).Value = $a
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
}

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("LP1912")
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws3 = $wb.Worksheets.Item("6203-6173")

# --- Sheet LP1912: header ---
$ws1.Cells.Item(2,1).Value = "Última actualización: 08:14:55"
$ws1.Cells.Item(3,1).Value = "Total filas: 68"

# --- Sheet LP1912: data rows ---
Set-RowValues $ws1 6 "05:47:29" "05:47" "17_ROMERO" 0 "LP1912"
Set-RowValues $ws1 7 "05:47:29" "06:09" "10_OLMOS" 22 "LP1912"
Set-RowValues $ws1 8 "06:15:23" "06:16" "215A_EL PATO" 1 "LP1912"
Set-RowValues $ws1 9 "06:15:23" "06:30" "23_HERNANDEZ" 15 "LP1912"
Set-RowValues $ws1 10 "06:15:23" "06:34" "11_ETCHEVERRY" 19 "LP1912"
Set-RowValues $ws1 11 "06:15:23" "06:39" "17X38_ROMERO" 24 "LP1912"
Set-RowValues $ws1 12 "06:15:23" "06:41" "16_SANTA ANA" 26 "LP1912"
Set-RowValues $ws1 13 "06:46:40" "06:56" "215A_EL PATO" 10 "LP1912"
Set-RowValues $ws1 14 "06:15:23" "06:57" "215A_EL PATO" 42 "LP1912"
Set-RowValues $ws1 15 "06:58:58" "06:58" "215A_EL PATO" 0 "LP1912"
Set-RowValues $ws1 16 "06:58:58" "06:58" "225_GOMEZ" 0 "LP1912"
Set-RowValues $ws1 17 "06:46:40" "06:59" "225_GOMEZ" 13 "LP1912"
Set-RowValues $ws1 18 "06:58:58" "07:15" "215C_EL PATO" 17 "LP1912"
Set-RowValues $ws1 19 "06:15:23" "07:16" "215C_EL PATO" 61 "LP1912"
Set-RowValues $ws1 20 "06:58:58" "07:18" "14_ABASTO" 20 "LP1912"
Set-RowValues $ws1 21 "06:46:40" "07:19" "14_ABASTO" 33 "LP1912"
Set-RowValues $ws1 22 "06:58:58" "07:20" "16_SANTA ANA" 22 "LP1912"
Set-RowValues $ws1 23 "06:46:40" "07:21" "23_HERNANDEZ" 35 "LP1912"
Set-RowValues $ws1 24 "06:15:23" "07:21" "16_SANTA ANA" 66 "LP1912"
Set-RowValues $ws1 25 "06:58:58" "07:22" "23_HERNANDEZ" 24 "LP1912"
Set-RowValues $ws1 26 "07:26:49" "07:29" "17X38_ROMERO" 3 "LP1912"
Set-RowValues $ws1 27 "06:58:58" "07:34" "10_OLMOS" 36 "LP1912"
Set-RowValues $ws1 28 "07:26:49" "07:35" "10_OLMOS" 9 "LP1912"
Set-RowValues $ws1 29 "07:26:49" "07:36" "27_EL RETIRO" 10 "LP1912"
Set-RowValues $ws1 30 "06:15:23" "07:37" "27_EL RETIRO" 82 "LP1912"
Set-RowValues $ws1 31 "07:26:49" "07:43" "215A_EL PATO" 17 "LP1912"
Set-RowValues $ws1 32 "06:58:58" "07:54" "14_ABASTO" 56 "LP1912"
Set-RowValues $ws1 33 "07:51:40" "07:55" "14_ABASTO" 4 "LP1912"
Set-RowValues $ws1 34 "06:58:58" "07:59" "17_ROMERO" 61 "LP1912"
Set-RowValues $ws1 35 "07:51:40" "08:00" "17_ROMERO" 9 "LP1912"
Set-RowValues $ws1 36 "07:26:49" "08:00" "16_SANTA ANA" 34 "LP1912"
Set-RowValues $ws1 37 "07:51:40" "08:01" "16_SANTA ANA" 10 "LP1912"
Set-RowValues $ws1 38 "07:51:40" "08:06" "23_HERNANDEZ" 15 "LP1912"
Set-RowValues $ws1 39 "07:51:40" "08:11" "10_OLMOS" 20 "LP1912"
Set-RowValues $ws1 40 "06:58:58" "08:12" "15X38_ABASTO" 74 "LP1912"
Set-RowValues $ws1 41 "07:51:40" "08:13" "15X38_ABASTO" 22 "LP1912"
Set-RowValues $ws1 42 "08:14:55" "08:14" "15X38_ABASTO" 0 "LP1912"
Set-RowValues $ws1 43 "06:58:58" "08:28" "15_ABASTO" 90 "LP1912"
Set-RowValues $ws1 44 "08:14:55" "08:28" "11_ETCHEVERRY" 14 "LP1912"
Set-RowValues $ws1 45 "07:51:40" "08:29" "11_ETCHEVERRY" 38 "LP1912"
Set-RowValues $ws1 46 "08:14:55" "08:29" "15_ABASTO" 15 "LP1912"
Set-RowValues $ws1 47 "06:58:58" "08:40" "16_P MOR-SANTA ANA" 102 "LP1912"
Set-RowValues $ws1 48 "08:14:55" "08:41" "16_P MOR-SANTA ANA" 27 "LP1912"
Set-RowValues $ws1 49 "08:14:55" "08:43" "215C_EL PATO" 29 "LP1912"
Set-RowValues $ws1 50 "07:51:40" "08:45" "23_HERNANDEZ" 54 "LP1912"
Set-RowValues $ws1 51 "08:14:55" "08:51" "23_HERNANDEZ" 37 "LP1912"
Set-RowValues $ws1 52 "06:58:58" "08:52" "23_HERNANDEZ" 114 "LP1912"
Set-RowValues $ws1 53 "08:14:55" "08:53" "215B_EL PATO" 39 "LP1912"
Set-RowValues $ws1 54 "08:14:55" "08:57" "215A_EL PATO" 43 "LP1912"
Set-RowValues $ws1 55 "07:51:40" "08:58" "215A_EL PATO" 67 "LP1912"
Set-RowValues $ws1 56 "08:14:55" "09:04" "10_OLMOS" 50 "LP1912"
Set-RowValues $ws1 57 "08:14:55" "09:06" "16_SANTA ANA" 52 "LP1912"
Set-RowValues $ws1 58 "07:26:49" "09:16" "27_EL RETIRO" 110 "LP1912"
Set-RowValues $ws1 59 "08:14:55" "09:17" "14_ABASTO" 63 "LP1912"
Set-RowValues $ws1 60 "08:14:55" "09:17" "27_EL RETIRO" 63 "LP1912"
Set-RowValues $ws1 61 "08:14:55" "09:18" "15X38_ABASTO" 64 "LP1912"
Set-RowValues $ws1 62 "07:51:40" "09:21" "27_EL RETIRO" 90 "LP1912"
Set-RowValues $ws1 63 "08:14:55" "09:28" "10_OLMOS" 74 "LP1912"
Set-RowValues $ws1 64 "07:51:40" "09:29" "10_OLMOS" 98 "LP1912"
Set-RowValues $ws1 65 "08:14:55" "09:31" "16_SANTA ANA" 77 "LP1912"
Set-RowValues $ws1 66 "08:14:55" "09:35" "23_HERNANDEZ" 81 "LP1912"
Set-RowValues $ws1 67 "08:14:55" "09:39" "15_ABASTO" 85 "LP1912"
Set-RowValues $ws1 68 "08:14:55" "09:41" "11_ETCHEVERRY" 87 "LP1912"
Set-RowValues $ws1 69 "08:14:55" "09:43" "16_P MOR-SANTA ANA" 89 "LP1912"
Set-RowValues $ws1 70 "08:14:55" "09:53" "10_OLMOS" 99 "LP1912"
Set-RowValues $ws1 71 "08:14:55" "09:58" "215C_EL PATO" 104 "LP1912"
Set-RowValues $ws1 72 "08:14:55" "10:05" "14_ABASTO" 111 "LP1912"
Set-RowValues $ws1 73 "08:14:55" "10:13" "17X38_ROMERO" 119 "LP1912"

# --- Sheet LP1912-215: header ---
$ws2.Cells.Item(2,1).Value = "Última actualización: 08:14:55"
$ws2.Cells.Item(3,1).Value = "Total filas: 12"

# --- Sheet LP1912-215: data rows ---
Set-RowValues $ws2 6 "06:15:23" "06:16" "215A_EL PATO" 1 "LP1912"
Set-RowValues $ws2 7 "06:46:40" "06:56" "215A_EL PATO" 10 "LP1912"
Set-RowValues $ws2 8 "06:15:23" "06:57" "215A_EL PATO" 42 "LP1912"
Set-RowValues $ws2 9 "06:58:58" "06:58" "215A_EL PATO" 0 "LP1912"
Set-RowValues $ws2 10 "06:58:58" "07:15" "215C_EL PATO" 17 "LP1912"
Set-RowValues $ws2 11 "06:15:23" "07:16" "215C_EL PATO" 61 "LP1912"
Set-RowValues $ws2 12 "07:26:49" "07:43" "215A_EL PATO" 17 "LP1912"
Set-RowValues $ws2 13 "08:14:55" "08:43" "215C_EL PATO" 29 "LP1912"
Set-RowValues $ws2 14 "08:14:55" "08:53" "215B_EL PATO" 39 "LP1912"
Set-RowValues $ws2 15 "08:14:55" "08:57" "215A_EL PATO" 43 "LP1912"
Set-RowValues $ws2 16 "07:51:40" "08:58" "215A_EL PATO" 67 "LP1912"
Set-RowValues $ws2 17 "08:14:55" "09:58" "215C_EL PATO" 104 "LP1912"

# --- Sheet 6203-6173: header ---
$ws3.Cells.Item(2,1).Value = "Última actualización: 08:14:55"
$ws3.Cells.Item(3,1).Value = "Total filas: 6"

# --- Sheet 6203-6173: data rows ---
Set-RowValues $ws3 6 "07:26:49" "07:42" "215A_LA PLATA" 16 "L6173"
Set-RowValues $ws3 7 "06:15:23" "07:43" "215A_LA PLATA" 88 "L6173"
Set-RowValues $ws3 8 "08:14:55" "08:35" "215A_LA PLATA" 21 "L6173"
Set-RowValues $ws3 9 "08:14:55" "08:50" "215C_LA PLATA" 36 "L6203"
Set-RowValues $ws3 10 "08:14:55" "09:20" "215A_LA PLATA" 66 "L6173"
Set-RowValues $ws3 11 "08:14:55" "10:12" "215C_LA PLATA" 118 "L6203"

